$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-09-09 06:35:16"
}
